$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'67.676.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.40%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.427.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.31%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'547.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.03%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'158.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.99%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.494"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.19%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'2.427.39"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.26%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = "'  -8.87%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  -1.68%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "'  -6.25%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = "'  -4.05%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'2.871.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.40%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'67.663.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.27%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("E16").Value = "'  -6.68%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'22.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.55%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'2.425.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.48%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'10.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.81%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'335.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.45%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'  -5.85%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'3.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.00%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.01%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "'  -5.31%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'65.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -5.10%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "'2.553.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.31%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").Value = "'3.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -7.34%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  +0.33%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'7.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -8.25%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'0.0₃0795"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -8.39%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'  -8.86%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.07%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'417.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -5.13%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'1.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.48%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = "'  -6.26%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'156.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.08%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'18.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.44%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  +0.03%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = "'  -4.80%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'17.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.44%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'  -5.33%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'4.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -6.83%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").Value = "'1.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.22%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -10.52%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'131.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.51%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'1.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -7.50%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'  -4.73%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "'  -2.59%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  -7.88%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.549"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.47%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.0896"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.50%  "
$ws.Range("E51").Style = "Normal"
